$updates = @{
    "G2" = [double]"3.685507"
    "H2" = [double]"11.056521"
    "I2" = [double]"0.3585631737883472"
    "J2" = [double]"0.3585631737883472"
    "K2" = [double]"2"
    "L2" = [double]"0.6666666666666666"
    "M2" = [double]"0.01848533333333334"
    "N2" = [double]"0.05545600000000001"
    "O2" = [double]"0.001625201930372746"
    "P2" = [double]"0.001625201930372746"
    "Q2" = [double]"0.06812782539733334"
    "R2" = [double]"0.6131504285760001"
    "S2" = [double]"0.0005827375622014003"
    "T2" = [double]"0.0005827375622014002"
    "G3" = [double]"3.685507"
    "H3" = [double]"11.056521"
    "I3" = [double]"0.3585631737883472"
    "J3" = [double]"0.3585631737883472"
    "O3" = [double]"0.002698334581238102"
    "P3" = [double]"0.002698334581238102"
    "Q3" = [double]"0.1131131238393333"
    "R3" = [double]"1.018018114554"
    "S3" = [double]"0.0009675234113915847"
    "T3" = [double]"0.0009675234113915846"
    "G4" = [double]"3.685507"
    "H4" = [double]"11.056521"
    "I4" = [double]"0.3585631737883472"
    "J4" = [double]"0.3585631737883472"
    "M4" = [double]"11.32499966666667"
    "N4" = [double]"33.974999"
    "O4" = [double]"0.9956764634883892"
    "P4" = [double]"0.995676463488389"
    "Q4" = [double]"41.73836554649768"
    "R4" = [double]"375.645289918479"
    "S4" = [double]"0.3570129128147543"
    "T4" = [double]"0.3570129128147542"
    "I5" = [double]"0.009647184430711629"
    "J5" = [double]"0.009647184430711629"
    "K5" = [double]"2"
    "L5" = [double]"0.6666666666666666"
    "M5" = [double]"0.01848533333333334"
    "N5" = [double]"0.05545600000000001"
    "O5" = [double]"0.001625201930372746"
    "P5" = [double]"0.001625201930372746"
    "Q5" = [double]"0.001832987168"
    "R5" = [double]"0.016496884512"
    "S5" = [double]"1.567862275945444E-05"
    "T5" = [double]"1.567862275945444E-05"
    "I6" = [double]"0.009647184430711629"
    "J6" = [double]"0.009647184430711629"
    "O6" = [double]"0.002698334581238102"
    "P6" = [double]"0.002698334581238102"
    "S6" = [double]"2.6031331360971E-05"
    "T6" = [double]"2.6031331360971E-05"
    "I7" = [double]"0.009647184430711629"
    "J7" = [double]"0.009647184430711629"
    "M7" = [double]"11.32499966666667"
    "N7" = [double]"33.974999"
    "O7" = [double]"0.9956764634883892"
    "P7" = [double]"0.995676463488389"
    "Q7" = [double]"1.122975641947"
    "R7" = [double]"10.106780777523"
    "S7" = [double]"0.009605474476591205"
    "T7" = [double]"0.009605474476591202"
    "G8" = [double]"6.493877"
    "H8" = [double]"19.481631"
    "I8" = [double]"0.6317896417809412"
    "J8" = [double]"0.6317896417809411"
    "K8" = [double]"2"
    "L8" = [double]"0.6666666666666666"
    "M8" = [double]"0.01848533333333334"
    "N8" = [double]"0.05545600000000001"
    "O8" = [double]"0.001625201930372746"
    "P8" = [double]"0.001625201930372746"
    "Q8" = [double]"0.1200414809706667"
    "R8" = [double]"1.080373328736"
    "S8" = [double]"0.001026785745411891"
    "T8" = [double]"0.001026785745411891"
    "G9" = [double]"6.493877"
    "H9" = [double]"19.481631"
    "I9" = [double]"0.6317896417809412"
    "J9" = [double]"0.6317896417809411"
    "O9" = [double]"0.002698334581238102"
    "P9" = [double]"0.002698334581238102"
    "Q9" = [double]"0.1993057436326667"
    "R9" = [double]"1.793751692694"
    "S9" = [double]"0.001704779838485546"
    "T9" = [double]"0.001704779838485546"
    "G10" = [double]"6.493877"
    "H10" = [double]"19.481631"
    "I10" = [double]"0.6317896417809412"
    "J10" = [double]"0.6317896417809411"
    "M10" = [double]"11.32499966666667"
    "N10" = [double]"33.974999"
    "O10" = [double]"0.9956764634883892"
    "P10" = [double]"0.995676463488389"
    "Q10" = [double]"73.54315486037434"
    "R10" = [double]"661.8883937433691"
    "S10" = [double]"0.6290580761970438"
    "T10" = [double]"0.6290580761970436"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}

Write-Output "Updated $($updates.Keys.Count) cells"
